$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JulyRaw")

# Header row
$ws.Cells.Item(1, 1).Value = "Library"
$ws.Cells.Item(1, 2).Value = "Items owned by this library checked out at this library this month"
$ws.Cells.Item(1, 3).Value = "Items owned by other libraries checked out at this library this month"
$ws.Cells.Item(1, 4).Value = "Total circulation this month"

# Data rows: row number, library name, B (owned/here), C (other/here), D (total)
$rows = @(
    @(2, "Atchison Public Library", 4953, 1540, 6493),
    @(3, "Baldwin City Public Library", 3420, 750, 4170),
    @(4, "Basehor Community Library", 9941, 1230, 11171),
    @(5, "Bern Community Library", 175, 46, 221),
    @(6, "Bonner Springs City Library", 6398, 1074, 7472),
    @(7, "Burlingame Community Library", 356, 331, 687),
    @(8, "Carbondale City Library", 677, 105, 782),
    @(9, "Centralia Community Library", 295, 82, 377),
    @(10, "Corning City Library", 75, $null, 75),
    @(11, "Digital Content", $null, $null, $null),
    @(12, "Doniphan County Library - Elwood", 107, 28, 135),
    @(13, "Doniphan County Library - Highland", 217, 89, 306),
    @(14, "Doniphan County Library - Troy", 535, 191, 726),
    @(15, "Doniphan County Library - Wathena", 303, 80, 383),
    @(16, "Effingham Community Library", 563, 81, 644),
    @(17, "Eudora Community Library", 1685, 653, 2338),
    @(18, "Everest, Barnes Reading Room", 90, 25, 115),
    @(19, "Hiawatha, Morrill Public Library", 2009, 465, 2474),
    @(20, "Highland Community College", 14, 12, 26),
    @(21, "Holton, Beck-Bookman Library", 1965, 498, 2463),
    @(22, "Horton Public Library", 406, 75, 481),
    @(23, "Lansing Community Library", 2533, 634, 3167),
    @(24, "Leavenworth Public Library", 11143, 1709, 12852),
    @(25, "Linwood Community Library", 631, 169, 800),
    @(26, "Louisburg Library", $null, $null, $null),
    @(27, "Lyndon Carnegie Library", 544, 219, 763),
    @(28, "McLouth Public Library", 304, 119, 423),
    @(29, "Meriden-Ozawkie Public Library", 2342, 524, 2866),
    @(30, "Northeast Kansas Library System", 17, 38, 55),
    @(31, "Nortonville Public Library", 343, 74, 417),
    @(32, "Osage City Library", 1512, 466, 1978),
    @(33, "Osawatomie Public Library", 1169, 477, 1646),
    @(34, "Oskaloosa Public Library", 683, 232, 915),
    @(35, "Ottawa Library", 8579, 1121, 9700),
    @(36, "Overbrook Public Library", 812, 182, 994),
    @(37, "Paola Free Library", 3895, 557, 4452),
    @(38, "Perry-Lecompton Community Library", 347, 35, 382),
    @(39, "Pomona Community Library", 127, 78, 205),
    @(40, "Prairie Hills Schools - Axtell Public School", $null, $null, $null),
    @(41, "Prairie Hills Schools - Sabetha Elementary School", $null, $null, $null),
    @(42, "Prairie Hills Schools - Sabetha High School", $null, $null, $null),
    @(43, "Prairie Hills Schools - Sabetha Middle School", $null, $null, $null),
    @(44, "Richmond Public Library", 576, 84, 660),
    @(45, "Rossville Community Library", 1687, 449, 2136),
    @(46, "Sabetha, Mary Cotton Library", 3881, 1010, 4891),
    @(47, "Seneca Free Library", 1958, 194, 2152),
    @(48, "Silver Lake Library", 948, 510, 1458),
    @(49, "Tonganoxie Public Library", 3777, 982, 4759),
    @(50, "Valley Falls, Delaware Township Library", 336, 130, 466),
    @(51, "Wellsville City Library", 939, 360, 1299),
    @(52, "Wetmore Public Library", 129, 96, 225),
    @(53, "Williamsburg Community Library", 182, 37, 219),
    @(54, "Winchester Public Library", 334, 393, 727)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    if ($row[2] -ne $null) { $ws.Cells.Item($r, 2).Value = $row[2] }
    if ($row[3] -ne $null) { $ws.Cells.Item($r, 3).Value = $row[3] }
    if ($row[4] -ne $null) { $ws.Cells.Item($r, 4).Value = $row[4] }
}

Write-Host "JulyRaw populated"
